{"js": "// Insert the \"\u03a6.15.1/\" protocol-prefix literal between the static\n// \".: \" label text and the \"${protocol}\" template placeholder in the\n// \"\u0391\u03c1. \u03a0\u03c1\u03c9\u03c4.:\" paragraph (top table, first cell).\n//\n// Before: \"\u0391\u03c1. \u03a0\u03c1\u03c9\u03c4.: ${protocol}\"\n// After:  \"\u0391\u03c1. \u03a0\u03c1\u03c9\u03c4.: \u03a6.15.1/${protocol}\"\n\nconst body = context.document.body;\n\n// \".: ${\" only occurs once in this template (right before the\n// \"protocol\" placeholder name), so it is a safe, unique anchor.\nconst results = body.search(\".: ${\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items,text\");\nawait context.sync();\n\nif (results.items.length !== 1) {\n  throw new Error(\n    `Expected exactly one match for \".: \\${\" but found ${results.items.length}`\n  );\n}\n\nconst target = results.items[0];\ntarget.insertText(\".: \u03a6.15.1/${\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Insert the \"\u03a6.15.1/\" protocol-prefix literal between the static\n# \".: \" label text and the \"${protocol}\" template placeholder in the\n# \"\u0391\u03c1. \u03a0\u03c1\u03c9\u03c4.:\" paragraph (top table, first cell).\n#\n# Before: \"\u0391\u03c1. \u03a0\u03c1\u03c9\u03c4.: ${protocol}\"\n# After:  \"\u0391\u03c1. \u03a0\u03c1\u03c9\u03c4.: \u03a6.15.1/${protocol}\"\n\n$d = $word.ActiveDocument\n\n# \".: ${\" only occurs once in this template (right before the\n# \"protocol\" placeholder name), so it is a safe, unique anchor.\n$rng = $d.Content\n$rng.Find.MatchWildcards = $false\n$found = $rng.Find.Execute(\".: `${\")\n\nif (-not $found) {\n    throw \"Could not find the '.: `${' anchor text\"\n}\n\n$rng.Text = \".: \u03a6.15.1/`${\"\n"}
